$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- The bus-scraper re-ran roughly one month later: refresh departure_date (C) ---
# --- and return_date (D) columns with the newly scraped serial dates.        ---
$ws.Range("C2").Value = 45108
$ws.Range("D2").Value = 45123
$ws.Range("C3").Value = 45108
$ws.Range("C4").Value = 45108
$ws.Range("C5").Value = 45108
$ws.Range("D5").Value = 45123
$ws.Range("C6").Value = 45108
$ws.Range("D6").Value = 45123
$ws.Range("D6").NumberFormat = "dd/mm/yy"
$ws.Range("C7").Value = 45108
$ws.Range("C8").Value = 45108
$ws.Range("C9").Value = 45108
$ws.Range("D9").Value = 45123
$ws.Range("C10").Value = 45108
$ws.Range("C11").Value = 45108
$ws.Range("D11").Value = 45123
$ws.Range("C12").Value = 45108
$ws.Range("C13").Value = 45108
$ws.Range("C14").Value = 45109
$ws.Range("C15").Value = 45109
$ws.Range("C16").Value = 45110
$ws.Range("C17").Value = 45110
$ws.Range("C18").Value = 45111
$ws.Range("D18").Value = 45126
$ws.Range("C19").Value = 45111
$ws.Range("C20").Value = 45112
$ws.Range("C21").Value = 45113
$ws.Range("C22").Value = 45114
$ws.Range("C23").Value = 45115
$ws.Range("D23").Value = 45130
$ws.Range("C24").Value = 45116
$ws.Range("C25").Value = 45117
$ws.Range("C26").Value = 45118
$ws.Range("D26").Value = 45133
$ws.Range("C27").Value = 45119
$ws.Range("C28").Value = 45120
$ws.Range("C29").Value = 45121
$ws.Range("D29").Value = 45136
$ws.Range("D29").NumberFormat = "dd/mm/yy"
$ws.Range("C30").Value = 45122

# --- New scraper run also appends two helper/debug columns (I, J) used by  ---
# --- the new buser_scrapper / clickbus debugger to stage upcoming queries. ---
# --- They are written out as blank, date-formatted cells for rows 2-30.   ---
$ws.Range("J2:J30").NumberFormat = "dd/mm/yy"
$ws.Range("I2").NumberFormat = "dd/mm/yy"
$ws.Range("I3:I30").NumberFormat = "dd/mm/yy"
$ws.Range("I3:I30").HorizontalAlignment = -4152

# --- Move the live selection to reflect where the scraper script left off. ---
$ws.Range("H21").Select()
